# Apply the FHIR IG publisher "deploy" update to the StructureDefinition
# spreadsheet: bump version/date, change publisher/jurisdiction metadata,
# and refresh the root-extension Short/Definition text on the Elements sheet.

$wb = $excel.ActiveWorkbook

# ---- Metadata sheet ----
$meta = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date updated
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty -> "Alvearie Team"
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 was a duplicated "Contact" / "No display for ContactDetail" row;
# turn it into "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was the second (duplicate) "Contact" row - remove it entirely so
# everything below shifts up by one row (A1:B21 -> A1:B20)
$meta.Rows.Item(11).Delete()

# ---- Elements sheet ----
$elements = $wb.Worksheets.Item("Elements")

# Root Extension row: Short / Definition refreshed to reference this profile
$elements.Range("K2").Value = "Problem Clinical Course"
$elements.Range("L2").Value = "Problem clinical course code"
